$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

function Set-PlainValue($ws, $cellRef, $val) {
    $ws.Range($cellRef).Value = $val
}

# Row 2
Set-PlainValue $ws 'D2' '64.339.48'
Set-PlainValue $ws 'E2' '  -3.18%  '

# Row 3
Set-PlainValue $ws 'D3' '3.174.46'
Set-PlainValue $ws 'E3' '  -4.39%  '

# Row 4
Set-TextValue $ws 'D4' '1.00'
Set-PlainValue $ws 'E4' '  +0.02%  '

# Row 5
Set-TextValue $ws 'D5' '569.61'
Set-PlainValue $ws 'E5' '  -2.86%  '

# Row 6
Set-TextValue $ws 'D6' '169.04'
Set-PlainValue $ws 'E6' '  -7.59%  '

# Row 7
Set-TextValue $ws 'D7' '0.610'
Set-PlainValue $ws 'E7' '  -5.63%  '

# Row 8
Set-PlainValue $ws 'E8' '  +0.06%  '

# Row 9
Set-PlainValue $ws 'D9' '3.175.01'
Set-PlainValue $ws 'E9' '  -4.29%  '

# Row 10
Set-TextValue $ws 'D10' '0.120'
Set-PlainValue $ws 'E10' '  -4.74%  '

# Row 11
Set-TextValue $ws 'D11' '6.72'
Set-PlainValue $ws 'E11' '  -0.92%  '

# Row 12
Set-TextValue $ws 'D12' '0.385'
Set-PlainValue $ws 'E12' '  -4.44%  '

# Row 13
Set-PlainValue $ws 'D13' '3.729.07'
Set-PlainValue $ws 'E13' '  -4.32%  '

# Row 14
Set-PlainValue $ws 'E14' '  -1.88%  '

# Row 15
Set-PlainValue $ws 'D15' '64.405.45'
Set-PlainValue $ws 'E15' '  -3.09%  '

# Row 16
Set-TextValue $ws 'D16' '25.39'
Set-PlainValue $ws 'E16' '  -3.96%  '

# Row 17
Set-PlainValue $ws 'E17' '  -3.60%  '

# Row 18
Set-PlainValue $ws 'D18' '3.181.06'
Set-PlainValue $ws 'E18' '  -4.13%  '

# Row 19
Set-TextValue $ws 'D19' '419.91'
Set-PlainValue $ws 'E19' '  -2.54%  '

# Row 20
Set-PlainValue $ws 'B20' 'Polkadot'
Set-PlainValue $ws 'C20' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws 'D20' '5.35'
Set-PlainValue $ws 'E20' '  -3.45%  '

# Row 21
Set-PlainValue $ws 'B21' 'Chainlink'
Set-PlainValue $ws 'C21' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws 'D21' '12.81'
Set-PlainValue $ws 'E21' '  -3.99%  '

# Row 22
Set-TextValue $ws 'D22' '7.03'
Set-PlainValue $ws 'E22' '  -5.61%  '

# Row 23
Set-PlainValue $ws 'E23' '  -0.30%  '

# Row 24
Set-TextValue $ws 'D24' '70.06'
Set-PlainValue $ws 'E24' '  -3.05%  '

# Row 25
Set-PlainValue $ws 'E25' '  +2.00%  '

# Row 26
Set-TextValue $ws 'D26' '0.488'
Set-PlainValue $ws 'E26' '  -5.69%  '

# Row 27
Set-PlainValue $ws 'E27' '  -7.24%  '

# Row 28
Set-TextValue $ws 'D28' '8.73'
Set-PlainValue $ws 'E28' '  -3.30%  '

# Row 29
Set-TextValue $ws 'D29' '0.999'
Set-PlainValue $ws 'E29' '  +0.02%  '

# Row 30
Set-TextValue $ws 'D30' '21.82'
Set-PlainValue $ws 'E30' '  -2.91%  '

# Row 31
Set-PlainValue $ws 'E31' '  -6.38%  '

# Row 33
Set-PlainValue $ws 'E33' '  -4.37%  '

# Row 34
Set-TextValue $ws 'D34' '6.32'
Set-PlainValue $ws 'E34' '  -4.76%  '

# Row 35
Set-TextValue $ws 'D35' '1.13'
Set-PlainValue $ws 'E35' '  -6.04%  '

# Row 36
Set-TextValue $ws 'D36' '157.19'
Set-PlainValue $ws 'E36' '  -1.44%  '

# Row 37
Set-PlainValue $ws 'E37' '  -6.75%  '

# Row 38
Set-PlainValue $ws 'D38' '2.696.89'
Set-PlainValue $ws 'E38' '  -6.59%  '

# Row 39
Set-PlainValue $ws 'E39' '  -7.79%  '

# Row 40
Set-TextValue $ws 'D40' '24.53'
Set-PlainValue $ws 'E40' '  -8.69%  '

# Row 41
Set-PlainValue $ws 'E41' '  -4.31%  '

# Row 42
Set-TextValue $ws 'D42' '39.08'
Set-PlainValue $ws 'E42' '  -3.13%  '

# Row 43
Set-TextValue $ws 'D43' '0.710'
Set-PlainValue $ws 'E43' '  -7.81%  '

# Row 44
Set-TextValue $ws 'D44' '5.67'
Set-PlainValue $ws 'E44' '  -6.13%  '

# Row 45
Set-TextValue $ws 'D45' '0.0620'
Set-PlainValue $ws 'E45' '  -7.22%  '

# Row 46
Set-TextValue $ws 'D46' '0.0261'
Set-PlainValue $ws 'E46' '  -4.00%  '

# Row 47
Set-TextValue $ws 'D47' '21.71'
Set-PlainValue $ws 'E47' '  -7.43%  '

# Row 48
Set-TextValue $ws 'D48' '293.98'
Set-PlainValue $ws 'E48' '  -7.63%  '

# Row 49
Set-PlainValue $ws 'B49' 'FirstDigitalUSD'
Set-PlainValue $ws 'C49' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws 'D49' '1.00'
Set-PlainValue $ws 'E49' '  -0.01%  '

# Row 50
Set-PlainValue $ws 'B50' 'dogwifhat'
Set-PlainValue $ws 'C50' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws 'D50' '2.00'
Set-PlainValue $ws 'E50' '  -13.80%  '

# Row 51
Set-TextValue $ws 'D51' '0.0991'
Set-PlainValue $ws 'E51' '  -5.20%  '
